# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-52) used to be listed newest-first
# (2009 .. 1709). They now need to read oldest-first (1709 .. 2009),
# and the "Valor Mora" (F) / "Salario Basico" (G) amounts for each of
# those periods are refreshed with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 16; Periodo = "1709"; Valor = 29509; Salario = 781242 },
    @{ Row = 17; Periodo = "1710"; Valor = 29509; Salario = 781242 },
    @{ Row = 18; Periodo = "1711"; Valor = 29509; Salario = 781242 },
    @{ Row = 19; Periodo = "1712"; Valor = 29509; Salario = 781242 },
    @{ Row = 20; Periodo = "1801"; Valor = 29509; Salario = 781242 },
    @{ Row = 21; Periodo = "1802"; Valor = 29509; Salario = 781242 },
    @{ Row = 22; Periodo = "1803"; Valor = 29509; Salario = 781242 },
    @{ Row = 23; Periodo = "1804"; Valor = 29509; Salario = 781242 },
    @{ Row = 24; Periodo = "1805"; Valor = 29509; Salario = 781242 },
    @{ Row = 25; Periodo = "1806"; Valor = 29509; Salario = 781242 },
    @{ Row = 26; Periodo = "1807"; Valor = 29509; Salario = 781242 },
    @{ Row = 27; Periodo = "1808"; Valor = 29509; Salario = 781242 },
    @{ Row = 28; Periodo = "1809"; Valor = 31249; Salario = 781242 },
    @{ Row = 29; Periodo = "1810"; Valor = 31249; Salario = 781242 },
    @{ Row = 30; Periodo = "1811"; Valor = 31249; Salario = 781242 },
    @{ Row = 31; Periodo = "1812"; Valor = 31249; Salario = 781242 },
    @{ Row = 32; Periodo = "1901"; Valor = 31249; Salario = 781242 },
    @{ Row = 33; Periodo = "1902"; Valor = 31249; Salario = 781242 },
    @{ Row = 34; Periodo = "1903"; Valor = 31249; Salario = 781242 },
    @{ Row = 35; Periodo = "1904"; Valor = 31249; Salario = 781242 },
    @{ Row = 36; Periodo = "1905"; Valor = 31249; Salario = 781242 },
    @{ Row = 37; Periodo = "1906"; Valor = 31249; Salario = 781242 },
    @{ Row = 38; Periodo = "1907"; Valor = 31249; Salario = 781242 },
    @{ Row = 39; Periodo = "1908"; Valor = 31249; Salario = 781242 },
    @{ Row = 40; Periodo = "1909"; Valor = 31249; Salario = 781242 },
    @{ Row = 41; Periodo = "1910"; Valor = 31249; Salario = 781242 },
    @{ Row = 42; Periodo = "1911"; Valor = 31249; Salario = 781242 },
    @{ Row = 43; Periodo = "1912"; Valor = 31249; Salario = 781242 },
    @{ Row = 44; Periodo = "2001"; Valor = 31249; Salario = 781242 },
    @{ Row = 45; Periodo = "2002"; Valor = 31249; Salario = 781242 },
    @{ Row = 46; Periodo = "2003"; Valor = 31249; Salario = 781242 },
    @{ Row = 47; Periodo = "2004"; Valor = 31249; Salario = 781242 },
    @{ Row = 48; Periodo = "2005"; Valor = 31249; Salario = 781242 },
    @{ Row = 49; Periodo = "2006"; Valor = 31249; Salario = 781242 },
    @{ Row = 50; Periodo = "2007"; Valor = 31249; Salario = 781242 },
    @{ Row = 51; Periodo = "2008"; Valor = 31249; Salario = 781242 },
    @{ Row = 52; Periodo = "2009"; Valor = 29166; Salario = 781242 }
)

foreach ($d in $data) {
    $ws.Range("E$($d.Row)").Value = $d.Periodo
    $ws.Range("F$($d.Row)").Value = $d.Valor
    $ws.Range("G$($d.Row)").Value = $d.Salario
}
